$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, shifting existing rows 12:96 down to 13:97.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new weekly record (same categorical
# fields as the row that is now 13, but a new date and new price figures).
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44685
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100112030
$ws.Range("G12").Value = "Poroto granado"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 440
$ws.Range("K12").Value = 26000
$ws.Range("L12").Value = 27000
$ws.Range("M12").Value = 26500
$ws.Range("N12").Value = "$/malla 25 kilos"
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 1060
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
